# edit.ps1 - apply "Corrections to problem 2" change:
#   - 8 brown socks -> 6 brown socks (with the _GoBack bookmark now sitting
#     right after the "6")
#   - rewrite the "three pairs ... collected." sentence, splitting it into
#     two runs
#   - add a brand-new "Fact: All other colors ..." paragraph right after it
#   - fill in the paragraph that used to hold the _GoBack bookmark (under
#     "Solving A:") with the new "18 socks would have to be drawn ..." text

$d = $word.ActiveDocument

function New-RunSplit($marker, $searchText) {
    # Forces a run boundary immediately BEFORE $searchText by temporarily
    # bookmarking that point and deleting the bookmark again; the run
    # split survives the bookmark's removal.
    $fr = $d.Content
    $fr.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $pt = $d.Range($fr.Start, $fr.Start)
    $d.Bookmarks.Add($marker, $pt)
    $d.Bookmarks($marker).Delete()
}

# ---------------------------------------------------------------------
# Change 1: "8 brown socks" -> "6" + _GoBack bookmark + " brown socks"
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("8 brown socks", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$numStart = $r.Start
$numRange = $d.Range($numStart, $numStart + 1)
$numRange.Text = "6"
$bmRange = $d.Range($numStart + 1, $numStart + 1)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------------
# Change 2a: rewrite the "three pairs ... collected." sentence
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute(" three pairs of each color a person would have to pair all but two socks in order to guarantee that the least common denominator (white) was collected. ", $true, $false, $false, $false, $false, $true, 1, $false, " three pairs of each color a person would have to solve for the least common denominator (white), because it has the least chance for collection.", 2)

# split the sentence into two runs, matching the source diff
New-RunSplit "TmpSplitA" "east common denominator (white), because it has the least chance for collection."

# ---------------------------------------------------------------------
# Change 2b: append a brand-new paragraph with the "Fact: All other ..." text
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("east common denominator (white), because it has the least chance for collection.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $r.Paragraphs(1)
$pEnd = $para.Range.End
$para.Range.InsertParagraphAfter()

$newPara = $d.Range($pEnd, $pEnd)
$newPara.Text = "Fact: All other colors would have to be collected in total (10 black) and (6 brown) so that the least common denominator could be solved."

New-RunSplit "TmpSplitB" "6 brown) so that the least common denominator could be solved."
New-RunSplit "TmpSplitC" " brown) so that the least common denominator could be solved."
New-RunSplit "TmpSplitD" "so that the least common denominator could be solved."

# ---------------------------------------------------------------------
# Change 3: fill the paragraph that used to contain the _GoBack bookmark
# (under "Solving A:") with the new summary sentence
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Solving A: ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$solvingPara = $r.Paragraphs(1)
$targetPara = $solvingPara.Next
$tr = $targetPara.Range
$tr.Text = "18 socks would have to be drawn in order to guarantee three pairs of each color of socks were matched."
